$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: split "mu" into "mu 1"/"mu 2", renumber Yxs columns, add L1/M1 ---
$ws.Range("H1").Value = "mu 1 [1/h]"
$ws.Range("I1").Value = "mu 2 [1/h]"
$ws.Range("J1").Value = "Yxs 1.1 [gx/gs]"
$ws.Range("K1").Value = "Yxs 1.2 [gx/gs]"
$ws.Range("L1").Value = "Yxs 2.1 [gx/gs]"
$ws.Range("M1").Value = "Yxs 2.2 [gx/gs]"

# New L1/M1 header cells need the same bold/border/center style as the other headers (copy from K1)
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null

# --- Data rows: shift Yxs values into new split columns, add mu2/Yxs2 values ---

# Row 2
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""

# Row 3
$ws.Range("I3").Value = 0.3346417012349824
$ws.Range("L3").Value = -14.0151515227274
$ws.Range("M3").Value = 14.0151515227274

# Row 4
$ws.Range("I4").Value = 0.8127218063500202
$ws.Range("J4").Value = -38.54545454545568
$ws.Range("K4").Value = 38.54545454545568
$ws.Range("L4").Value = 10.51515149999989
$ws.Range("M4").Value = 10.51515149999989

# Row 5
$ws.Range("I5").Value = 0.4243727592997257
$ws.Range("J5").Value = -3.299688472897195
$ws.Range("K5").Value = 3.299688472897195
$ws.Range("L5").Value = -1.788174138401557
$ws.Range("M5").Value = 1.788174138401557

# Row 6
$ws.Range("I6").Value = 0.2746536833672744
$ws.Range("J6").Value = -2.568350168686871
$ws.Range("K6").Value = 2.568350168686871
$ws.Range("L6").Value = -1.708424909890113
$ws.Range("M6").Value = 1.708424909890113

# Row 7
$ws.Range("I7").Value = 0.1442387591429934
$ws.Range("J7").Value = -1.73448275862069
$ws.Range("K7").Value = 1.73448275862069
$ws.Range("L7").Value = -0.633777777333333
$ws.Range("M7").Value = 0.633777777333333

# Row 8
$ws.Range("I8").Value = 0.2217057364726398
$ws.Range("J8").Value = -1.684195525010689
$ws.Range("K8").Value = 1.684195525010689
$ws.Range("L8").Value = -1.538119086811352
$ws.Range("M8").Value = 1.538119086811352

# Row 9
$ws.Range("I9").Value = 0.07735908689312881
$ws.Range("J9").Value = -1.340620445658963
$ws.Range("K9").Value = 1.340620445658963
$ws.Range("L9").Value = -0.409424488991888
$ws.Range("M9").Value = 0.409424488991888

# Row 10
$ws.Range("I10").Value = 0.07007962478060897
$ws.Range("J10").Value = -1.274203850807555
$ws.Range("K10").Value = 1.274203850807555
$ws.Range("L10").Value = -0.8026607538802666
$ws.Range("M10").Value = 0.8026607538802666

# Row 11
$ws.Range("I11").Value = 0.03719083035979676
$ws.Range("J11").Value = -1.162319534099067
$ws.Range("K11").Value = 1.162319534099067
$ws.Range("L11").Value = -0.3852978441064636
$ws.Range("M11").Value = 0.3852978441064636

# Row 12
$ws.Range("I12").Value = 0.03619552273593946
$ws.Range("J12").Value = 0.3442443842321475
$ws.Range("K12").Value = 0.3442443842321475
$ws.Range("L12").Value = 0.1115874013155464
$ws.Range("M12").Value = 0.1115874013155464

# Row 13
$ws.Range("I13").Value = 0.03801145000663797
$ws.Range("J13").Value = 0.2886938913263584
$ws.Range("K13").Value = 0.2886938913263584
$ws.Range("L13").Value = 0.1003260225251926
$ws.Range("M13").Value = 0.1003260225251926

# Row 14
$ws.Range("I14").Value = -0.06954913092679149
$ws.Range("J14").Value = 0.2305471595849802
$ws.Range("K14").Value = 0.2305471595849802
$ws.Range("L14").Value = -0.1300404631645038
$ws.Range("M14").Value = 0.1300404631645038

# Row 15
$ws.Range("I15").Value = -0.01970310092027705
$ws.Range("J15").Value = 0.1802583930214906
$ws.Range("K15").Value = 0.1802583930214906
$ws.Range("L15").Value = -0.04301075277419363
$ws.Range("M15").Value = 0.04301075277419363

# Row 16
$ws.Range("I16").Value = -0.006451635241295772
$ws.Range("J16").Value = 0.1596327558122316
$ws.Range("K16").Value = 0.1596327558122316
$ws.Range("L16").Value = -0.01043092763543897
$ws.Range("M16").Value = 0.01043092763543897

# Row 17
$ws.Range("I17").Value = 0.1810127751752835
$ws.Range("J17").Value = 0.1650565328603363
$ws.Range("K17").Value = 0.1650565328603363
$ws.Range("L17").Value = 0.2264813926689029
$ws.Range("M17").Value = 0.2264813926689029

# Row 18
$ws.Range("I18").Value = -0.02883454103336998
$ws.Range("J18").Value = 0.1537079804978188
$ws.Range("K18").Value = 0.1537079804978188
$ws.Range("L18").Value = -0.03383659910054636
$ws.Range("M18").Value = 0.03383659910054636

# Row 19
$ws.Range("I19").Value = 0.1444724872797511
$ws.Range("J19").Value = 0.187037036984127
$ws.Range("K19").Value = 0.187037036984127
$ws.Range("L19").Value = 0.4024322826344469
$ws.Range("M19").Value = 0.4024322826344469

# Row 20
$ws.Range("I20").Value = 0.02448929697916169
$ws.Range("J20").Value = 0.1848466320697034
$ws.Range("K20").Value = 0.1848466320697034
$ws.Range("L20").Value = 0.127823694214876
$ws.Range("M20").Value = 0.127823694214876

# Row 21
$ws.Range("I21").Value = 0.0387149826119924
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = ""
